# Refresh the Victoria "Key outbreaks" cluster list: new clusters added,
# several resolved/renamed clusters removed, and active-case counts updated
# for rows 2 through 34 (header row 1 is untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2; Name = '139 Highett St Apartment Complex Richmond'; Cases = 10 },
    @{ Row = 3; Name = '3175 The Bays Aged Care Facility Hastings'; Cases = 14 },
    @{ Row = 4; Name = '3563 Embracia Aged Care Reservoir'; Cases = 22 },
    @{ Row = 5; Name = 'Apartment Complex Fawkner'; Cases = 10 },
    @{ Row = 6; Name = 'Australian Lamb Colac East'; Cases = 13 },
    @{ Row = 7; Name = 'Bread Solutions Braeside'; Cases = 13 },
    @{ Row = 8; Name = 'CS Square Caroline Springs'; Cases = 13 },
    @{ Row = 9; Name = 'Carton Finishing Pty. Ltd. Campbellfield'; Cases = 12 },
    @{ Row = 10; Name = 'Cedar Meats Australia Brooklyn'; Cases = 10 },
    @{ Row = 11; Name = 'Community Kids Bayswater Early Education Centre Bayswater North'; Cases = 18 },
    @{ Row = 12; Name = 'Costco Wholesale Epping'; Cases = 13 },
    @{ Row = 13; Name = 'Ermha365 Residential Disability Care Services Doveton'; Cases = 10 },
    @{ Row = 14; Name = 'FedEx Station Melbourne Airport'; Cases = 14 },
    @{ Row = 15; Name = 'Green Leaves Early Learning Centre Highlands Craigieburn'; Cases = 14 },
    @{ Row = 16; Name = 'Guardian Childcare Caulfield'; Cases = 14 },
    @{ Row = 17; Name = 'Kool Kidz Childcare Narre Warren'; Cases = 16 },
    @{ Row = 18; Name = 'Lantmannen Unibake Australia Mordialloc'; Cases = 26 },
    @{ Row = 19; Name = 'Melbourne Assessment Prison West Melbourne'; Cases = 10 },
    @{ Row = 20; Name = 'MyCentre Childcare Broadmeadows'; Cases = 10 },
    @{ Row = 21; Name = 'Nido Early School Ascot Vale'; Cases = 29 },
    @{ Row = 22; Name = 'Nido Early School Glenroy'; Cases = 24 },
    @{ Row = 23; Name = 'Northern Health Northern Hospital Epping Emergency Department Tier 1B'; Cases = 52 },
    @{ Row = 24; Name = 'Northern Health The Northern Hospital Epping'; Cases = 21 },
    @{ Row = 25; Name = 'Social Gathering Warrnambool 28 September'; Cases = 17 },
    @{ Row = 26; Name = 'St Margaret''s Primary School OSHC Maribyrnong'; Cases = 11 },
    @{ Row = 27; Name = 'St Vincents Hospital Emergency Department Melbourne'; Cases = 35 },
    @{ Row = 28; Name = 'The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B'; Cases = 16 },
    @{ Row = 29; Name = 'The Royal Talbot Rehabilitation Centre'; Cases = 12 },
    @{ Row = 30; Name = 'Visy Recycling Springvale'; Cases = 31 },
    @{ Row = 31; Name = 'Wallaby Childcare Wollert'; Cases = 16 },
    @{ Row = 32; Name = 'Werribee Mercy Hospital Emergency Department'; Cases = 20 },
    @{ Row = 33; Name = 'Western Health Footscray Hospital Emergency Department'; Cases = 10 },
    @{ Row = 34; Name = 'Western Health Sunshine Hospital Emergency Department'; Cases = 16 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Name
    $ws.Cells.Item($r.Row, 2).Value = $r.Cases
}

Write-Host "Updated $($rows.Count) rows"
